$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.678.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.062.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.87%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.36"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +8.09%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.055.81"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.75%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.00%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.21"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.94"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.568.52"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.697.01"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.065.48"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.64"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.55"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.63%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.97"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.32%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +10.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.05"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.45"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.77"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.99"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "468.12"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.167.94"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0814"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.60%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.61%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "28.42"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +14.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.82%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.44%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0510"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.97"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.09%  "
